$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Done" column header in D1
$ws.Range("D1").Value = "Done"

# Fill D2:D69 with "y" for every existing data row
$ws.Range("D2:D69").Value = "y"

# Move active selection to D70 (just below the new data), matching the
# recorded end-state of the edit
$ws.Range("D70").Select()
